$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values are stored as text cells in the source sheet, so force text
# number format before assignment to avoid Excel auto-converting the
# numeric-looking / percent-looking strings into numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-5.68%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.29%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.096"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.72%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07840"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.91%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.316"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.96%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.689"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-12.89%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9199"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.08%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1062"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.48%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1758"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.40%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09073"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.96%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.28%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.182"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-16.43%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1057"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.17%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001296"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.21%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005829"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.83%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.379"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.26%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.567"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.70%"
$ws.Range("E18").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1387"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.18%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.85%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04171"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.30%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001220"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.05%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004137"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.15%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.76%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.45%"
$ws.Range("E26").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02427"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-9.90%"
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05261"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.79%"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.07%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1352"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.07%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007483"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.37%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002022"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.64%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008069"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.54%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3363"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.22%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006770"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.10%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.42%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003199"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-8.33%"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004124"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.67%"
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002112"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.42%"
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002011"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.42%"
$ws.Range("E51").Style = "Normal"
